# Support importing of doubles entries
#
# This populates the "Club" column (D) on the "entries" sheet for the
# doubles pairs that share a club/region (CHN, EUR, KOR), and moves the
# active-tab/selection from the "players" sheet back to "entries" (cell
# D10), matching the reviewer's last position after adding the data.

$wb = $excel.ActiveWorkbook
$entries = $wb.Worksheets.Item("entries")

# --- Fill in the Club column for the doubles teams that share one ---
$entries.Range("D2").Value = "CHN"
$entries.Range("D3").Value = "CHN"
$entries.Range("D6").Value = "EUR"
$entries.Range("D7").Value = "EUR"
$entries.Range("D8").Value = "KOR"

# --- Move the active sheet/selection back onto "entries" ---
$entries.Activate() | Out-Null
$entries.Range("D10").Select() | Out-Null
